# Insert a new "year" column before the existing "cpi" column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing column A ("cpi") to column B by inserting a new column at A.
$ws.Columns.Item(1).Insert()

# New column A header and values (years 1990-2013).
$ws.Range("A1").Value = "year"

$years = 1990..2013
for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]
}

# Update selection to match the target state.
$ws.Range("E12").Select()
